# Insert a new data row at row 39 (pushing existing rows 39:92 down to 40:93)
# and populate it with a new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44771
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = 100112026
$ws.Cells.Item(39, 7).Value = "Haba"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 90
$ws.Cells.Item(39, 11).Value = 18000
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = 18000
$ws.Cells.Item(39, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 720
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
